# feat: add 2022-Q3 data
#
# 1) "总计" sheet: insert a new row for 2022-Q1 (pushing the former row2/row3
#    values along) and put the new 2022-Q3 figures into row 2.
# 2) Duplicate the "2022-Q2" sheet (so it inherits the same header/cell
#    styles), rename the duplicate to "2022-Q3", drop its extra data row and
#    overwrite the remaining two rows with the 2022-Q3 fund holdings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Make room for a 4th row by cloning row 3's formatting into row 4, then
# fill it with what used to be the 2022-Q1 figures.
$summary.Range("A3:D3").Copy($summary.Range("A4:D4"))
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 3
$summary.Range("D4").Value = 0.95

# Row 3 now becomes the old 2022-Q2 figures.
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 1.39

# Row 2 becomes the new 2022-Q3 figures.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.53

# ---------------------------------------------------------------------
# 2. Create the "2022-Q3" sheet from a copy of "2022-Q2"
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The source sheet has 3 data rows; 2022-Q3 only needs 2, drop the 3rd.
$q3.Rows(4).Delete()

# Overwrite the remaining two data rows with the 2022-Q3 holdings.
$q3.Range("B2").Value = "007216"
$q3.Range("C2").Value = "浙商中华预期高股息C"
$q3.Range("D2").Value = "4.40"
$q3.Range("E2").Value = "88.55"
$q3.Range("F2").Value = "7.53"
$q3.Range("G2").Value = "0.3313"
$q3.Range("H2").Value = 3

$q3.Range("B3").Value = "007178"
$q3.Range("C3").Value = "浙商中华预期高股息A"
$q3.Range("D3").Value = "2.59"
$q3.Range("E3").Value = "88.55"
$q3.Range("F3").Value = "7.53"
$q3.Range("G3").Value = "0.1950"
$q3.Range("H3").Value = 3
